$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove all hyperlinks from the sheet (kept as plain text/urls in cells).
$ws.Hyperlinks.Delete()

# 2. Insert a new row after row 195 (becomes row 196) for the new
#    "虚拟贸易展会 (Virtual Trade Show)" entry under "虚拟活动（Virtual Events）".
$ws.Rows.Item(196).Insert()
$ws.Range("A196").Value = "虚拟活动（Virtual Events）"
$ws.Range("B196").Value = "虚拟贸易展会（Virtual Trade Show）"
$ws.Range("C196").Value = "盈拓展览｜https://vr.zhizhan360.com/STONETEC-2022/"
$ws.Range("D196").Value = 20220729

# 3. Append four new rows at the bottom of the sheet (rows 221-224): three
#    more video-player entries, plus a new Web3.0 row.
$ws.Range("A221").Value = "更多"
$ws.Range("B221").Value = "视频播放器"
$ws.Range("C221").Value = "射手影音播放器｜https://github.com/chiflix/splayerx"
$ws.Range("D221").Value = 20220729

$ws.Range("A222").Value = "更多"
$ws.Range("B222").Value = "视频播放器"
$ws.Range("C222").Value = "electron-player｜https://github.com/c10342/player"
$ws.Range("D222").Value = 20220729

$ws.Range("A223").Value = "更多"
$ws.Range("B223").Value = "视频播放器"
$ws.Range("C223").Value = "SGPlayer｜https://github.com/libobjc/SGPlayer"
$ws.Range("D223").Value = 20220729

$ws.Range("A224").Value = "更多"
$ws.Range("B224").Value = "Web3.0"
$ws.Range("C224").Value = "Meson Network｜https://meson.network/"
$ws.Range("D224").Value = 20220729

# 4. Update the saved view: scroll down and select C205 (matches how the
#    author left the sheet positioned after the edit).
$ws.Range("C205").Select()
$excel.ActiveWindow.ScrollRow = 202
